# Artisan s7.xlsx help sheet — remove the term "slave" from the PID
# Control dialog description (commit: "removes the term 'slave' from the
# project where possible").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit -----------------------------------------------------
# Cell A8 holds the PID Control paragraph; drop "slave" (x2) from the text.
$ws.Range("A8").Value = "The PID Control dialog can operate a connected PID using the given PID registers to set the p-i-d parameters and the set value (SV). S7 commands can be specified to turn the PID on and off from that PID Control dialog. See the help page in the Events Dialog for documentation of available S7 write commands."

# --- Cosmetic state carried over from the resave ----------------------
# Row 4 (the wrapped "S7 SETTINGS" intro paragraph) reflows to a shorter
# auto height once resaved.
$ws.Rows.Item(4).RowHeight = 22.7

# Cursor/selection position stored at save time.
$ws.Range("A13").Select()

# Header/footer margins re-expressed at full precision (1.3 cm) by the
# resaving application.
$ws.PageSetup.HeaderMargin = 36.8503937007874
$ws.PageSetup.FooterMargin = 36.8503937007874
